$d = $word.ActiveDocument

# P6 (Objetivos body) text -> becomes old Programa-resumido text
$search0 = "Apresentar os principais conceitos sobre as transformações de fases em materiais metálicos, poliméricos e cerâmicos abrangendo transformações difusionais e não-difusionais, a conceituação sobre nucleação e crescimento (aspectos termodinâmicos e cinéticos) e sua relação com problemas práticos encontrados nas indústrias de processamento e de transformação de materiais."
$replace0 = "Difusão no estado sólido. Difusão em materiais não-metálicos. Recuperação, recristalização e crescimento de grão. Solidificação. Precipitação no estado sólido. Cinética de transformação no sistema Fe-C e em ligas não-ferrosas. Transformação de fases em vidros e cerâmicas. Transformação de fases em materiais poliméricos. Atividade experimental."
$rng = $d.Paragraphs.Item(6).Range
$ok = $rng.Find.Execute($search0, $true, $false, $false, $false, $false, $true, 0, $false, $replace0, 2)
Write-Output "Para 6: $ok"

# P8 run0 (Gilberto) -> Objetivos text (trailing br preserved)
$search1 = "5009972 - Gilberto Carvalho Coelho"
$replace1 = "Apresentar os principais conceitos sobre as transformações de fases em materiais metálicos, poliméricos e cerâmicos abrangendo transformações difusionais e não-difusionais, a conceituação sobre nucleação e crescimento (aspectos termodinâmicos e cinéticos) e sua relação com problemas práticos encontrados nas indústrias de processamento e de transformação de materiais."
$rng = $d.Paragraphs.Item(8).Range
$ok = $rng.Find.Execute($search1, $true, $false, $false, $false, $false, $true, 0, $false, $replace1, 2)
Write-Output "Para 8: $ok"

# P8 run1 (Hugo) -> Programa text1 + br + text2 (no trailing br)
$search2 = "984972 - Hugo Ricardo Zschommler Sandim"
$replace2 = "Introdução à difusão no estado sólido. Coeficiente de difusão. Leis de Fick. Difusão em soluções diluídas e na presença de um gradiente de concentração. Efeito Kirkendall. Apresentar os fundamentos teóricos pertinentes à transformação de fases em materiais metálicos, cerâmicos e poliméricos. Apresentar os conceitos fundamentais associados à nucleação (homogênea e heterogênea), ao crescimento e à cinética de transformação de fases. Aspectos microestruturais relevantes em fundidos. Precipitação no estado sólido. Descrição das principais transformações de fase no estado sólido no sistema Fe-C e em algumas ligas não-ferrosas. Curvas TTT e CCT (TRC). " + [char]11 + "Realização de atividade experimental (8 horas-aula) versando sobre tópicos da ementa para consolidação dos conhecimentos teóricos. Viagem Didática complementar."
$rng = $d.Paragraphs.Item(8).Range
$ok = $rng.Find.Execute($search2, $true, $false, $false, $false, $false, $true, 0, $false, $replace2, 2)
Write-Output "Para 8: $ok"

# P10 run0 (Programa resumido text) -> Metodo text1+br+text2+br+text3 (no trailing br)
$search3 = "Difusão no estado sólido. Difusão em materiais não-metálicos. Recuperação, recristalização e crescimento de grão. Solidificação. Precipitação no estado sólido. Cinética de transformação no sistema Fe-C e em ligas não-ferrosas. Transformação de fases em vidros e cerâmicas. Transformação de fases em materiais poliméricos. Atividade experimental."
$replace3 = "Esta é uma disciplina de caráter fundamental, exigindo dedicação individual para assimilação das definições e conceitos. Isto envolve leitura concentrada para fixação dos conceitos teóricos e realização de exercícios numéricos. " + [char]11 + "O aluno será avaliado ao longo do semestre por duas avaliações escritas (P1 e P2) correspondendo a 80% do total da nota final, em pesos iguais, e uma atividade experimental (AE) correspondendo a 20% da nota final. Um relatório circunstanciado sobre o experimento atribuído, além da apresentação oral dos resultados, integra a avaliação da atividade experimental (8 horas-aula). " + [char]11 + "O desenvolvimento do aluno ao longo do curso será aferido e estimulado por meio de discussões sobre um dado tema, porém sem a atribuição de nota, por conta da subjetividade envolvida."
$rng = $d.Paragraphs.Item(10).Range
$ok = $rng.Find.Execute($search3, $true, $false, $false, $false, $false, $true, 0, $false, $replace3, 2)
Write-Output "Para 10: $ok"

# P12 run0 (Programa text1+br+text2) -> Critério text (NotaFinal)
$search4 = "Introdução à difusão no estado sólido. Coeficiente de difusão. Leis de Fick. Difusão em soluções diluídas e na presença de um gradiente de concentração. Efeito Kirkendall. Apresentar os fundamentos teóricos pertinentes à transformação de fases em materiais metálicos, cerâmicos e poliméricos. Apresentar os conceitos fundamentais associados à nucleação (homogênea e heterogênea), ao crescimento e à cinética de transformação de fases. Aspectos microestruturais relevantes em fundidos. Precipitação no estado sólido. Descrição das principais transformações de fase no estado sólido no sistema Fe-C e em algumas ligas não-ferrosas. Curvas TTT e CCT (TRC). " + [char]11 + "Realização de atividade experimental (8 horas-aula) versando sobre tópicos da ementa para consolidação dos conhecimentos teóricos. Viagem Didática complementar."
$replace4 = "A Nota final (NF) será calculada da seguinte maneira: NF = (0,4P1 + 0,4P2 + 0,2AE)"
$rng = $d.Paragraphs.Item(12).Range
$ok = $rng.Find.Execute($search4, $true, $false, $false, $false, $false, $true, 0, $false, $replace4, 2)
Write-Output "Para 12: $ok"

# P14 run5 (Norma text1+br+text2, no trailing br) -> Gilberto
$search5 = "Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). " + [char]11 + "Média final = (NF + PR) / 2"
$replace5 = "5009972 - Gilberto Carvalho Coelho"
$rng = $d.Paragraphs.Item(14).Range
$ok = $rng.Find.Execute($search5, $true, $false, $false, $false, $false, $true, 0, $false, $replace5, 2)
Write-Output "Para 14: $ok"

# P14 run1 (Metodo text1+br+text2+br+text3, trailing br preserved) -> Norma text1+br+text2
$search6 = "Esta é uma disciplina de caráter fundamental, exigindo dedicação individual para assimilação das definições e conceitos. Isto envolve leitura concentrada para fixação dos conceitos teóricos e realização de exercícios numéricos. " + [char]11 + "O aluno será avaliado ao longo do semestre por duas avaliações escritas (P1 e P2) correspondendo a 80% do total da nota final, em pesos iguais, e uma atividade experimental (AE) correspondendo a 20% da nota final. Um relatório circunstanciado sobre o experimento atribuído, além da apresentação oral dos resultados, integra a avaliação da atividade experimental (8 horas-aula). " + [char]11 + "O desenvolvimento do aluno ao longo do curso será aferido e estimulado por meio de discussões sobre um dado tema, porém sem a atribuição de nota, por conta da subjetividade envolvida."
$replace6 = "Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). " + [char]11 + "Média final = (NF + PR) / 2"
$rng = $d.Paragraphs.Item(14).Range
$ok = $rng.Find.Execute($search6, $true, $false, $false, $false, $false, $true, 0, $false, $replace6, 2)
Write-Output "Para 14: $ok"

# P14 run3 (CriterioTxt NotaFinal, trailing br preserved) -> Bibliography list (10 items, 9 internal br)
$search7 = "A Nota final (NF) será calculada da seguinte maneira: NF = (0,4P1 + 0,4P2 + 0,2AE)"
$replace7 = "1. CALLISTER Jr, W.D., RETHWISCH, D.G. Ciência e Engenharia de Materiais: Uma Introdução, 8ª ed., LTC Editora, 2013." + [char]11 + "2. ASKELAND, D.R., PHULÉ, P.P., Ciência e Engenharia dos Materiais, CENGAGE, São Paulo, 2008." + [char]11 + "3. SHACKELFORD, J.F., Ciência dos Materiais, 6a. ed., Pearson, 2008." + [char]11 + "4. GARCIA, A. Solidificação: Fundamentos e Aplicações. Ed. UNICAMP, 2001." + [char]11 + "5. READEY, D. W. Kinetics in Materials Science and Engineering. CRC Press, 1st. Ed. 2016." + [char]11 + "6. SHEWMON, P.G. Diffusion in solids. McGraw-Hill, 1963. " + [char]11 + "7. SHEWMON, P.G. Phase transformation in metals. McGraw-Hill, 1969. " + [char]11 + "8. HUMPHREYS, F.J, HATHERLY, M. Recrystallization and related annealing phenomena. Pergamon, 2004. " + [char]11 + "9. BILLMEYER JR., F.W. Textbook of Polymer Science. John Wiley & Sons, New York, 1984. " + [char]11 + "10. WILSON, E.A. Worked examples in the kinetics and thermodynamics of phase transformations. CRC Press, 1a. Ed., 1981"
$rng = $d.Paragraphs.Item(14).Range
$ok = $rng.Find.Execute($search7, $true, $false, $false, $false, $false, $true, 0, $false, $replace7, 2)
Write-Output "Para 14: $ok"

# P16 run0 (Bibliography list) -> Hugo
$search8 = "1. CALLISTER Jr, W.D., RETHWISCH, D.G. Ciência e Engenharia de Materiais: Uma Introdução, 8ª ed., LTC Editora, 2013." + [char]11 + "2. ASKELAND, D.R., PHULÉ, P.P., Ciência e Engenharia dos Materiais, CENGAGE, São Paulo, 2008." + [char]11 + "3. SHACKELFORD, J.F., Ciência dos Materiais, 6a. ed., Pearson, 2008." + [char]11 + "4. GARCIA, A. Solidificação: Fundamentos e Aplicações. Ed. UNICAMP, 2001." + [char]11 + "5. READEY, D. W. Kinetics in Materials Science and Engineering. CRC Press, 1st. Ed. 2016." + [char]11 + "6. SHEWMON, P.G. Diffusion in solids. McGraw-Hill, 1963. " + [char]11 + "7. SHEWMON, P.G. Phase transformation in metals. McGraw-Hill, 1969. " + [char]11 + "8. HUMPHREYS, F.J, HATHERLY, M. Recrystallization and related annealing phenomena. Pergamon, 2004. " + [char]11 + "9. BILLMEYER JR., F.W. Textbook of Polymer Science. John Wiley & Sons, New York, 1984. " + [char]11 + "10. WILSON, E.A. Worked examples in the kinetics and thermodynamics of phase transformations. CRC Press, 1a. Ed., 1981"
$replace8 = "984972 - Hugo Ricardo Zschommler Sandim"
$rng = $d.Paragraphs.Item(16).Range
$ok = $rng.Find.Execute($search8, $true, $false, $false, $false, $false, $true, 0, $false, $replace8, 2)
Write-Output "Para 16: $ok"
